$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -19.2
$ws.Range("E12").Value = "40.8/140"
